$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.475.33"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "'1.573.93"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("E7").Value = "  -0.95%  "
$ws.Range("D8").Value = "'49.89"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.3410"
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("D10").Value = "'1.150"
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("D11").Value = "'0.07558"
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "'21.28"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").Value = "'6.052"
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").Value = "'6.973"
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("D16").Value = "'1.575.86"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "'91.07"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").Value = "'0.06757"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").Value = "'6.308"
$ws.Range("E21").Value = "  +1.33%  "
$ws.Range("D22").Value = "'16.37"
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("D23").Value = "'12.16"
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("D24").Value = "'22.467.65"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("D26").Value = "'2.643"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'20.06"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("D28").Value = "'149.00"
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("D29").Value = "'5.059"
$ws.Range("E29").Value = "  +0.60%  "
$ws.Range("D30").Value = "'125.67"
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("D31").Value = "'1.749.40"
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").Value = "'1.079"
$ws.Range("E32").Value = "  +9.80%  "
$ws.Range("D33").Value = "'6.238"
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("D34").Value = "'2.011"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "'9.846"
$ws.Range("E35").Value = "  -2.60%  "
$ws.Range("D36").Value = "'0.08394"
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("D37").Value = "'0.02489"
$ws.Range("E37").Value = "  -1.94%  "
$ws.Range("D38").Value = "'0.2306"
$ws.Range("E38").Value = "  -0.35%  "
$ws.Range("D39").Value = "'0.06558"
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("D40").Value = "'1.336"
$ws.Range("E40").Value = "  -2.41%  "
$ws.Range("D41").Value = "'5.469"
$ws.Range("E41").Value = "  +1.07%  "
$ws.Range("D42").Value = "'11.37"
$ws.Range("E42").Value = "  -0.57%  "
$ws.Range("D43").Value = "'0.6252"
$ws.Range("E43").Value = "  -1.99%  "
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "'14.01"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").Value = "'3.812"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").Value = "'0.5848"
$ws.Range("E47").Value = "  -2.00%  "
$ws.Range("D48").Value = "'130.22"
$ws.Range("E48").Value = "  +4.28%  "
$ws.Range("D49").Value = "'2.083"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("D50").Value = "'1.223"
$ws.Range("E50").Value = "  -5.33%  "
